# "Generate Report for Handback" - refresh the localization-status report
# for the 56cfe6a0-514c-4956-b444-8ae31df7596b handback: it now resolves
# to a non-latest version, so both locale sheets (zh-cn, de-de) get their
# row 5 (Latest Target File / Latest Handback File / Latest Handback
# DateTime / Error Detail) filled in, plus widened columns to fit the new
# long text.

$wb = $excel.ActiveWorkbook

$errMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45e79554b7a4ea35914c8f3a591c861b81d7f889/e2e/56cfe6a0-514c-4956-b444-8ae31df7596b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc4d6374f343adff5d6073ae784b275a12294752/e2e/56cfe6a0-514c-4956-b444-8ae31df7596b.md."

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc4d6374f343adff5d6073ae784b275a12294752/e2e/56cfe6a0-514c-4956-b444-8ae31df7596b.md"

# Hyperlink font color used throughout this workbook for the "HyperLink"
# cell style (RGB 0x6495ED == VBA color value 15570276).
$hlColor = 15570276

$sheets = @(
    @{ Name = "zh-cn"; HandbackFile = "56cfe6a0-514c-4956-b444-8ae31df7596b.2c2c4daaf564ed21d42c4ed746e2068eca9f34e7.zh-cn.xlf"; HandbackDate = "2016-11-15 16:29:39" },
    @{ Name = "de-de"; HandbackFile = "56cfe6a0-514c-4956-b444-8ae31df7596b.2c2c4daaf564ed21d42c4ed746e2068eca9f34e7.de-de.xlf"; HandbackDate = "2016-11-15 16:29:58" }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the columns that now hold long file names / URLs.
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Row 5 is the 56cfe6a0-... handback, which the report now flags as
    # stale: fill in its target file, handback file, handback datetime
    # and the error detail explaining the mismatch.
    $ws.Hyperlinks.Add($ws.Range("I5"), $latestUrl, "", "", "56cfe6a0-514c-4956-b444-8ae31df7596b.md")
    $ws.Range("I5").Font.Underline = $true
    $ws.Range("I5").Font.Color = $hlColor

    $ws.Range("J5").Value = $info.HandbackFile
    $ws.Range("K5").Value = $info.HandbackDate
    $ws.Range("P5").Value = $errMsg
}
